$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18 (20:00) - Tuesday, Wednesday, Thursday, Friday
$ws.Range("C18").Value = "[Anderson-Tornearia-1NB, Anderson-Tornearia-1NB, Anderson-Tornearia-1NB, Anderson-Tornearia-1NB]"
$ws.Range("D18").Value = "Humberto-M.T.F.-"
$ws.Range("E18").Value = "[Sandro M.-Comandos Eletricos-1NB, Sandro M.-Comandos Eletricos-1NB, Sandro M.-Comandos Eletricos-1NB, Sandro M.-Comandos Eletricos-1NB]"
$ws.Range("F18").Value = "[Rachel-T.M; Metalicos-1NB, Rachel-T.M; Metalicos-1NB, Rachel-T.M; Metalicos-1NB, Rachel-T.M; Metalicos-1NB]"

# Row 19 (20:50) - Monday, Tuesday, Wednesday, Thursday, Friday
$ws.Range("B19").Value = "[Gisele-E. D. N. D.-1NB, Gisele-E. D. N. D.-1NB, Gisele-E. D. N. D.-1NB, Gisele-E. D. N. D.-1NB]"
$ws.Range("C19").Value = "[Mayra-T. NãoMetalicos-1NB, Mayra-T. NãoMetalicos-1NB, Mayra-T. NãoMetalicos-1NB, Mayra-T. NãoMetalicos-1NB]"
$ws.Range("D19").Value = "Humberto-M.T.F.-"
$ws.Range("E19").Value = "[Valmir-Caldeiraria-1NB, Valmir-Caldeiraria-1NB, Valmir-Caldeiraria-1NB, Valmir-Caldeiraria-1NB]"
$ws.Range("F19").Value = "[Suzanny-Trat. Termicos-1NB, Suzanny-Trat. Termicos-1NB, Suzanny-Trat. Termicos-1NB, Suzanny-Trat. Termicos-1NB]"

# Row 20 (21:40) - Monday, Tuesday, Wednesday, Thursday, Friday
$ws.Range("B20").Value = "[Humberto-Desenho tecnico mecanico-T2-1NB, Humberto-Desenho tecnico mecanico-T2-1NB, Humberto-Desenho tecnico mecanico-T2-1NB, Humberto-Desenho tecnico mecanico-T2-1NB]"
$ws.Range("C20").Value = "[Aselmo-Manut. Mot. End.-1NB, Aselmo-Manut. Mot. End.-1NB, Aselmo-Manut. Mot. End.-1NB, Aselmo-Manut. Mot. End.-1NB]"
$ws.Range("D20").Value = "[Weslei-Desenho tecnico mecanico-T1-1NB, Weslei-Desenho tecnico mecanico-T1-1NB, Weslei-Desenho tecnico mecanico-T1-1NB, Weslei-Desenho tecnico mecanico-T1-1NB]"
$ws.Range("E20").Value = "Gilberto-M.T.R.M.-"
$ws.Range("F20").Value = "[Victor S.-Ajustagem-1NB, Victor S.-Ajustagem-1NB, Victor S.-Ajustagem-1NB, Victor S.-Ajustagem-1NB]"

# Row 21 (22:35) - Monday, Tuesday, Wednesday, Thursday
$ws.Range("B21").Value = "[Weslei-Metrologia 1-1NB, Weslei-Metrologia 1-1NB, Weslei-Metrologia 1-1NB, Weslei-Metrologia 1-1NB]"
$ws.Range("C21").Value = "[Suzanny-Metalografia-1NB, Suzanny-Metalografia-1NB, Suzanny-Metalografia-1NB, Suzanny-Metalografia-1NB]"
$ws.Range("D21").Value = "[Andre B.-Elet. Digi. Básica-1NB, Andre B.-Elet. Digi. Básica-1NB, Andre B.-Elet. Digi. Básica-1NB, Andre B.-Elet. Digi. Básica-1NB]"
$ws.Range("E21").Value = "Gilberto-M.T.R.M.-"
